$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks numeric ("6", "120.00", ...) as genuine
# TEXT rather than letting Excel auto-convert it to a Number. We do this by
# entering a string-literal formula (="6") and then converting that cell to
# a static value in place via Copy + PasteSpecial(values) -- this keeps the
# cell's number format / style untouched (no quote-prefix style needed).
function Set-TextValue($range, [string]$val) {
    $escaped = $val -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Rename the worksheet tab to the player's name
$ws.Name = "Jaydev Unadkat"

# Insert a new "matchNo" column at A, shifting the existing columns right
$ws.Columns.Item(1).Insert()

# Header row
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# The pre-existing data row (currently row 2) belongs to the "7th" match.
# Record its matchNo now, before shifting it down to row 3.
$ws.Range("A2").Value = "7th"

# Insert a fresh row above it for the "54th" match (new row 2)
$ws.Rows.Item(2).Insert()

# Row 2: 54th match vs Kolkata Knight Riders
$ws.Range("A2").Value = "54th"
$ws.Range("B2").Value = "Rajasthan Royals"
$ws.Range("C2").Value = "Jaydev Unadkat"
$ws.Range("D2").Value = "c Shakib Al Hasan b Ferguson"
Set-TextValue $ws.Range("E2") "6"
Set-TextValue $ws.Range("F2") "5"
Set-TextValue $ws.Range("G2") "1"
Set-TextValue $ws.Range("H2") "0"
Set-TextValue $ws.Range("I2") "120.00"
$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Sharjah"
$ws.Range("L2").Value = "October 07"
$ws.Range("M2").Value = "KKR won by 86 runs"

# Row 3: the original "7th" match data, now shifted into the new layout
$ws.Range("B3").Value = "Rajasthan Royals"
$ws.Range("C3").Value = "Jaydev Unadkat"
# "states" is an empty string for this match (no dismissal text) -- write it
# as a genuine empty TEXT cell via a scratch-cell copy/paste-special trick
# (assigning "" directly would just clear the cell instead of leaving an
# empty-but-present text value).
$ws.Range("Z1").Formula = "="""""
$ws.Range("Z1").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
Set-TextValue $ws.Range("E3") "11"
Set-TextValue $ws.Range("F3") "7"
Set-TextValue $ws.Range("G3") "0"
Set-TextValue $ws.Range("H3") "1"
Set-TextValue $ws.Range("I3") "157.14"
$ws.Range("J3").Value = "Delhi Capitals"
$ws.Range("K3").Value = "Wankhede"
$ws.Range("L3").Value = "April 15"
$ws.Range("M3").Value = "Royals won by 3 wickets (with 2 balls remaining)"

# Row 4: new 12th match vs Chennai Super Kings
$ws.Range("A4").Value = "12th"
$ws.Range("B4").Value = "Rajasthan Royals"
$ws.Range("C4").Value = "Jaydev Unadkat"
$ws.Range("D4").Value = "c Jadeja b Thakur"
Set-TextValue $ws.Range("E4") "24"
Set-TextValue $ws.Range("F4") "17"
Set-TextValue $ws.Range("G4") "2"
Set-TextValue $ws.Range("H4") "1"
Set-TextValue $ws.Range("I4") "141.17"
$ws.Range("J4").Value = "Chennai Super Kings"
$ws.Range("K4").Value = "Wankhede"
$ws.Range("L4").Value = "April 19"
$ws.Range("M4").Value = "Super Kings won by 45 runs"
